$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (price / 1h volume change) and row 43/44 coin swap
# Column D values that look numeric are forced to remain plain text (matching the
# source data, which mixes thousand-separated / malformed numeric strings) by
# briefly applying a Text number format, then resetting the cell style to Normal
# so no lingering number-format style is left on the cell.
$ws.Range("D2").Value = '68.179.00'
$ws.Range("E2").Value = '  +1.20%  '
$ws.Range("D3").Value = '3.275.94'
$ws.Range("E3").Value = '  +0.83%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '588.42'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.98%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '185.68'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.47%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.135'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.56%  '
$ws.Range("E10").Value = '  -0.04%  '
$ws.Range("E11").Value = '  +0.94%  '
$ws.Range("D12").Value = '3.844.08'
$ws.Range("E12").Value = '  +0.86%  '
$ws.Range("E13").Value = '  +0.39%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.56'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.15%  '
$ws.Range("D15").Value = '68.162.93'
$ws.Range("E15").Value = '  +1.32%  '
$ws.Range("E16").Value = '  +2.42%  '
$ws.Range("D17").Value = '3.276.61'
$ws.Range("E17").Value = '  +0.87%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.87'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.04%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.68'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.36%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '381.56'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.44%  '
$ws.Range("E21").Value = '  +1.85%  '
$ws.Range("E22").Value = '  +0.00%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '71.48'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.44%  '
$ws.Range("E24").Value = '  +2.74%  '
$ws.Range("E25").Value = '  +0.89%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.191'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +6.72%  '
$ws.Range("E27").Value = '  -1.98%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.998'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.35%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.85'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.96%  '
$ws.Range("E30").Value = '  +1.15%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '22.95'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.87%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.17'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.12%  '
$ws.Range("E33").Value = '  +0.71%  '
$ws.Range("E34").Value = '  -0.01%  '
$ws.Range("E35").Value = '  +3.01%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '162.67'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.48%  '
$ws.Range("E37").Value = '  +0.47%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.839'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.07%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.77'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.01%  '
$ws.Range("E40").Value = '  +1.82%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '26.65'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.67%  '
$ws.Range("E42").Value = '  +4.61%  '
$ws.Range("B43").Value = 'Hedera'
$ws.Range("C43").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0693'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.95%  '
$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.39'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.28%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '25.50'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.05%  '
$ws.Range("D46").Value = '2.641.31'
$ws.Range("E46").Value = '  -4.21%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '343.10'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.48%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0285'
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '32.47'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.44%  '
$ws.Range("E50").Value = '  +0.74%  '
$ws.Range("E51").Value = '  -0.08%  '
